# "regenerate orders with updates distance/sizes"
#
# The underlying trial-order table encodes experiment parameters inside
# text tokens (e.g. "Face04_D51_S30", "Fixation_D51_l.png", the bare
# "D51"/"D64"/"D80" distance codes and the bare "S30" size code). This
# edit just renumbers those codes:
#   D51 -> D55
#   D64 -> D69
#   D80 -> D86
#   S30 -> S31
# (S25 and S20 are untouched.) Every place these tokens appear -
# Condition, Filename_Left, Filename_Right, Distance and Size columns -
# needs the same substring substitution, so driving it through
# Range.Replace over the whole used range covers every affected cell in
# one pass, regardless of row/column layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$used = $ws.UsedRange

# Range.Replace(What, Replacement, LookAt, SearchOrder, MatchCase, MatchByte, SearchFormat, ReplaceFormat)
# LookAt:=2 (xlPart) so the substring matches inside longer tokens like
# "Face04_D51_S30" / "Face04_D51_S30_l.png"; MatchCase:=$true so we only
# touch the intended uppercase D/S codes.
$xlPart = 2
$xlByRows = 1

$used.Replace("D51", "D55", $xlPart, $xlByRows, $true)
$used.Replace("D64", "D69", $xlPart, $xlByRows, $true)
$used.Replace("D80", "D86", $xlPart, $xlByRows, $true)
$used.Replace("S30", "S31", $xlPart, $xlByRows, $true)
